# 1st changes of mifos to finflux
#
# The "Repayment schedule" worksheet gets a new (blank) column inserted
# immediately before the existing "Late" column (column N), pushing the
# "Late" / "heading" / "Outstanding" columns one slot to the right.
# Afterwards the "Repayment schedule" tab becomes the active tab (it was
# "Transactions" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make this sheet active/selected (was "Transactions" before the edit).
$ws.Activate()

# Insert a new blank column before column N ("Late"), shifting
# N->O, O->P, P->Q and copying the formatting from the column
# immediately to the left (column M), matching the width Excel applies
# when inserting a column in the middle of a formatted table.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the selection to where the user left it after the insert.
$ws.Range("R5").Select() | Out-Null
